# Commit message: "Update Name of Algo"
# This applies the updated KNN imputation results to the result_data_KNN
# worksheet: a set of recomputed numeric cell values across columns A-E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.64
$ws.Range("C3").Value = -12.774
$ws.Range("E8").Value = 16.867
$ws.Range("D19").Value = -8.001999999999999
$ws.Range("A21").Value = -20.257
$ws.Range("A23").Value = -20.775
$ws.Range("E23").Value = 16.46
$ws.Range("C24").Value = -12.424
$ws.Range("D24").Value = -7.755
$ws.Range("A25").Value = -21.711
$ws.Range("E26").Value = 16.419
$ws.Range("B27").Value = 5.626
$ws.Range("D30").Value = -7.191
$ws.Range("B31").Value = 5.351
$ws.Range("D31").Value = -8.086
$ws.Range("D33").Value = -7.69
$ws.Range("E37").Value = 16.498
$ws.Range("B39").Value = 7.915000000000001
$ws.Range("B48").Value = 5.431999999999999
$ws.Range("E48").Value = 17.075
$ws.Range("B51").Value = 5.962000000000001
$ws.Range("B52").Value = 5.295
$ws.Range("A53").Value = -22.01
$ws.Range("B55").Value = 4.508
$ws.Range("D55").Value = -8.221
$ws.Range("B56").Value = 5.144
$ws.Range("A57").Value = -21.768
$ws.Range("B57").Value = 6.531000000000001
$ws.Range("C57").Value = -12.98
$ws.Range("A59").Value = -22.363
$ws.Range("C61").Value = -13.508
$ws.Range("E62").Value = 16.612
$ws.Range("D65").Value = -7.869
$ws.Range("E66").Value = 16.985
$ws.Range("A69").Value = -21.541
$ws.Range("C70").Value = -12.097
$ws.Range("D70").Value = -7.449
$ws.Range("B73").Value = 7.561
$ws.Range("D75").Value = -7.702
$ws.Range("A79").Value = -21.115
$ws.Range("A83").Value = -22.134
$ws.Range("D83").Value = -8.512
$ws.Range("C86").Value = -13.597
$ws.Range("B89").Value = 4.971
$ws.Range("E89").Value = 17.192
$ws.Range("B90").Value = 5.811
$ws.Range("A93").Value = -21.536
$ws.Range("E94").Value = 17.228
$ws.Range("D96").Value = -7.447
$ws.Range("D97").Value = -8.166
$ws.Range("C98").Value = -12.45
$ws.Range("C100").Value = -12.204
$ws.Range("C102").Value = -13.556
